# Apply the "4th june" commit edits to the Login sheet of LoginData.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: new URL, new username (Administrator), new password
$ws.Range("A2").Value = "https://lab.singtel.tetherfi.cloud:45443/OCMUI"
$ws.Range("B2").Value = "Administrator"
$ws.Range("C2").Value = "pSHS-Iq;DXfKp;dAw;Lfufub&CEL*-tD"

# Update the active selection to B12 (previously C12)
$ws.Range("B12").Select()
